$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 3
$ws.Range("L2").Value = 0.3
$ws.Range("N2").Value = 3.15
$ws.Range("O2").Value = "The evidence report is reasonably structured with clear sections and consistent formatting, earning a score of 4 in structure and formatting. However, the evidence extraction quality is only acceptable, as the report includes some fragmented snippets and lacks full sentence quotations, resulting in a score of 3. The coverage of representativeness dimensions is addressed but lacks depth and specificity, particularly in geographic and language distribution, leading to a score of 3. Relevance and faithfulness are maintained, but the report occasionally lacks direct grounding in the source documents, scoring a 3. Missing disclosures are identified, but not comprehensively, also scoring a 3. Overall, the report is useful for auditing but could be more detailed and precise, resulting in a score of 3 for audit usefulness."
